# Updates to latest 4.0
#
# 1) About sheet: remove the EU/US pre-tax biofuel cost adjustment block
#    (rows 26-28), shrinking the used range from A1:B28 to A1:B24.
# 2) ICtPSFfL sheet: row 7 ("biofuel diesel") formulas no longer multiply
#    by the now-removed About!$B$27/About!$B$28 ratio.

$wb = $excel.ActiveWorkbook

# --- 1) About: delete rows 26:28 ---------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A26:B28").EntireRow.Delete() | Out-Null

# --- 2) ICtPSFfL: strip the About!$B$27/About!$B$28 multiplier ---------
$wsFuel = $wb.Worksheets.Item("ICtPSFfL")

for ($col = 2; $col -le 37; $col++) {
    $colLetter = ($wsFuel.Cells.Item(35, $col).Address($false, $false)) -replace '35', ''
    $wsFuel.Cells.Item(7, $col).Formula = "=MAX(Calcs!" + $colLetter + "35,0)"
}
